$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new issue row pair (37:38) for the new issue, mirroring the
# layout used by the other two-row (merged) issue blocks. Do this first
# so the new shared string for this issue is appended to the shared
# string table before any other new strings introduced below.
# Match the same center-aligned formatting used by every other data row.
$ws.Range("A37:D38").HorizontalAlignment = -4108

$ws.Range("A37").Value2 = "Mutliple Users same device Location Persistence"
$ws.Range("B37").Value2 = "Open"
$ws.Range("C37").Value2 = "Tejas"
$ws.Range("D37").Value2 = "High"
$ws.Range("A38").Value2 = ""
$ws.Range("B38").Value2 = ""
$ws.Range("C38").Value2 = ""
$ws.Range("D38").Value2 = ""

# Merge the new row pair the same way as other issue rows.
$ws.Range("A37:A38").Merge() | Out-Null
$ws.Range("B37:B38").Merge() | Out-Null
$ws.Range("C37:C38").Merge() | Out-Null
$ws.Range("D37:D38").Merge() | Out-Null

# "Places API integration" row (row 13): mark Resolved by Tejas
$ws.Range("B13").Value2 = "Resolved"
$ws.Range("C13").Value2 = "Tejas"

# Rename "Nearby List" -> "Nearby Contacts List" (row 17)
$ws.Range("A17").Value2 = "Nearby Contacts List"

# "Firebase integration" row (row 35): priority High -> Medium
$ws.Range("D35").Value2 = "Medium"

# Update the sheet view: selection now points at the "Checkbox removal
# near contact" row (A19:A20), and the frozen/top-left cell was reset
# back to the default (no explicit topLeftCell override).
$ws.Range("A19:A20").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
